# Dynamic environment url implementation
# Code cleaning and optimization regarding driver object
#
# The sheet's C1 cell used to hold a hardcoded environment URL
# (https://d5e0000019ce6eai) rendered as a live hyperlink. Remove the
# hyperlink and its text so the cell goes back to being a plain, empty
# input cell (matching the blank cells used elsewhere on the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("C1")

# Drop the hyperlink relationship that lived on C1.
$target.Hyperlinks.Delete()

# Clear the old URL text (this also drops the now-unused shared string).
$target.Value = ""

# Re-base the cell's look on a neighbouring plain cell instead of the
# special blue/underlined "hyperlink" font that C1 used to carry.
$ws.Range("B2").Copy()
$target.PasteSpecial(-4122)
